$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Cells.Item(2,4).Value = "64.773.23"
$ws.Cells.Item(3,4).Value = "3.095.38"
$ws.Cells.Item(5,4).Formula = "=""559.80"""
$ws.Cells.Item(5,4).Copy()
$ws.Cells.Item(5,4).PasteSpecial(-4163)
$ws.Cells.Item(6,4).Formula = "=""144.45"""
$ws.Cells.Item(6,4).Copy()
$ws.Cells.Item(6,4).PasteSpecial(-4163)
$ws.Cells.Item(8,4).Value = "3.093.05"
$ws.Cells.Item(10,4).Formula = "=""7.04"""
$ws.Cells.Item(10,4).Copy()
$ws.Cells.Item(10,4).PasteSpecial(-4163)
$ws.Cells.Item(12,4).Formula = "=""0.466"""
$ws.Cells.Item(12,4).Copy()
$ws.Cells.Item(12,4).PasteSpecial(-4163)
$ws.Cells.Item(14,4).Formula = "=""35.35"""
$ws.Cells.Item(14,4).Copy()
$ws.Cells.Item(14,4).PasteSpecial(-4163)
$ws.Cells.Item(15,4).Value = "3.601.86"
$ws.Cells.Item(16,4).Value = "64.828.57"
$ws.Cells.Item(17,4).Value = "3.100.77"
$ws.Cells.Item(19,4).Formula = "=""6.83"""
$ws.Cells.Item(19,4).Copy()
$ws.Cells.Item(19,4).PasteSpecial(-4163)
$ws.Cells.Item(20,4).Formula = "=""480.87"""
$ws.Cells.Item(20,4).Copy()
$ws.Cells.Item(20,4).PasteSpecial(-4163)
$ws.Cells.Item(21,4).Formula = "=""13.78"""
$ws.Cells.Item(21,4).Copy()
$ws.Cells.Item(21,4).PasteSpecial(-4163)
$ws.Cells.Item(22,4).Formula = "=""0.679"""
$ws.Cells.Item(22,4).Copy()
$ws.Cells.Item(22,4).PasteSpecial(-4163)
$ws.Cells.Item(23,4).Formula = "=""7.48"""
$ws.Cells.Item(23,4).Copy()
$ws.Cells.Item(23,4).PasteSpecial(-4163)
$ws.Cells.Item(24,4).Formula = "=""13.41"""
$ws.Cells.Item(24,4).Copy()
$ws.Cells.Item(24,4).PasteSpecial(-4163)
$ws.Cells.Item(25,4).Formula = "=""81.08"""
$ws.Cells.Item(25,4).Copy()
$ws.Cells.Item(25,4).PasteSpecial(-4163)
$ws.Cells.Item(26,4).Formula = "=""1.00"""
$ws.Cells.Item(26,4).Copy()
$ws.Cells.Item(26,4).PasteSpecial(-4163)
$ws.Cells.Item(27,4).Formula = "=""2.78"""
$ws.Cells.Item(27,4).Copy()
$ws.Cells.Item(27,4).PasteSpecial(-4163)
$ws.Cells.Item(28,4).Formula = "=""8.17"""
$ws.Cells.Item(28,4).Copy()
$ws.Cells.Item(28,4).PasteSpecial(-4163)
$ws.Cells.Item(31,4).Formula = "=""26.05"""
$ws.Cells.Item(31,4).Copy()
$ws.Cells.Item(31,4).PasteSpecial(-4163)
$ws.Cells.Item(33,4).Formula = "=""2.47"""
$ws.Cells.Item(33,4).Copy()
$ws.Cells.Item(33,4).PasteSpecial(-4163)
$ws.Cells.Item(34,4).Formula = "=""5.66"""
$ws.Cells.Item(34,4).Copy()
$ws.Cells.Item(34,4).PasteSpecial(-4163)
$ws.Cells.Item(35,4).Formula = "=""6.22"""
$ws.Cells.Item(35,4).Copy()
$ws.Cells.Item(35,4).PasteSpecial(-4163)
$ws.Cells.Item(36,4).Formula = "=""55.08"""
$ws.Cells.Item(36,4).Copy()
$ws.Cells.Item(36,4).PasteSpecial(-4163)
$ws.Cells.Item(37,4).Formula = "=""469.45"""
$ws.Cells.Item(37,4).Copy()
$ws.Cells.Item(37,4).PasteSpecial(-4163)
$ws.Cells.Item(38,4).Formula = "=""0.0409"""
$ws.Cells.Item(38,4).Copy()
$ws.Cells.Item(38,4).PasteSpecial(-4163)
$ws.Cells.Item(39,4).Formula = "=""0.0828"""
$ws.Cells.Item(39,4).Copy()
$ws.Cells.Item(39,4).PasteSpecial(-4163)
$ws.Cells.Item(40,4).Formula = "=""2.92"""
$ws.Cells.Item(40,4).Copy()
$ws.Cells.Item(40,4).PasteSpecial(-4163)
$ws.Cells.Item(41,4).Value = "3.007.36"
$ws.Cells.Item(42,4).Formula = "=""8.27"""
$ws.Cells.Item(42,4).Copy()
$ws.Cells.Item(42,4).PasteSpecial(-4163)
$ws.Cells.Item(44,4).Formula = "=""28.27"""
$ws.Cells.Item(44,4).Copy()
$ws.Cells.Item(44,4).PasteSpecial(-4163)
$ws.Cells.Item(45,4).Formula = "=""0.257"""
$ws.Cells.Item(45,4).Copy()
$ws.Cells.Item(45,4).PasteSpecial(-4163)
$ws.Cells.Item(49,4).Value = "0.0₃0518"
$ws.Cells.Item(50,4).Formula = "=""117.71"""
$ws.Cells.Item(50,4).Copy()
$ws.Cells.Item(50,4).PasteSpecial(-4163)

# Update Volume(1h) (column E) values
$ws.Cells.Item(2,5).Value = "  +2.92%  "
$ws.Cells.Item(3,5).Value = "  +1.20%  "
$ws.Cells.Item(4,5).Value = "  +0.13%  "
$ws.Cells.Item(5,5).Value = "  +2.23%  "
$ws.Cells.Item(6,5).Value = "  +5.97%  "
$ws.Cells.Item(7,5).Value = "  +0.03%  "
$ws.Cells.Item(8,5).Value = "  +1.34%  "
$ws.Cells.Item(9,5).Value = "  +0.45%  "
$ws.Cells.Item(10,5).Value = "  +12.57%  "
$ws.Cells.Item(11,5).Value = "  +1.88%  "
$ws.Cells.Item(12,5).Value = "  +2.69%  "
$ws.Cells.Item(13,5).Value = "  +3.04%  "
$ws.Cells.Item(14,5).Value = "  +0.40%  "
$ws.Cells.Item(15,5).Value = "  +1.82%  "
$ws.Cells.Item(16,5).Value = "  +3.09%  "
$ws.Cells.Item(17,5).Value = "  +1.76%  "
$ws.Cells.Item(18,5).Value = "  -0.57%  "
$ws.Cells.Item(19,5).Value = "  +1.45%  "
$ws.Cells.Item(20,5).Value = "  -1.26%  "
$ws.Cells.Item(21,5).Value = "  +2.49%  "
$ws.Cells.Item(22,5).Value = "  -0.80%  "
$ws.Cells.Item(23,5).Value = "  +5.13%  "
$ws.Cells.Item(24,5).Value = "  +9.38%  "
$ws.Cells.Item(25,5).Value = "  -1.04%  "
$ws.Cells.Item(26,5).Value = "  -0.02%  "
$ws.Cells.Item(27,5).Value = "  +1.85%  "
$ws.Cells.Item(28,5).Value = "  +2.99%  "
$ws.Cells.Item(29,5).Value = "  +5.21%  "
$ws.Cells.Item(30,5).Value = "  +0.26%  "
$ws.Cells.Item(31,5).Value = "  -0.08%  "
$ws.Cells.Item(32,5).Value = "  +0.12%  "
$ws.Cells.Item(33,5).Value = "  +3.67%  "
$ws.Cells.Item(34,5).Value = "  -2.01%  "
$ws.Cells.Item(35,5).Value = "  +4.29%  "
$ws.Cells.Item(36,5).Value = "  -1.15%  "
$ws.Cells.Item(37,5).Value = "  +2.01%  "
$ws.Cells.Item(38,5).Value = "  +4.17%  "
$ws.Cells.Item(39,5).Value = "  +2.27%  "
$ws.Cells.Item(40,5).Value = "  +17.66%  "
$ws.Cells.Item(41,5).Value = "  -5.63%  "
$ws.Cells.Item(42,5).Value = "  +0.77%  "
$ws.Cells.Item(43,5).Value = "  -2.80%  "
$ws.Cells.Item(44,5).Value = "  +6.44%  "
$ws.Cells.Item(45,5).Value = "  +3.30%  "
$ws.Cells.Item(47,5).Value = "  +5.20%  "
$ws.Cells.Item(48,5).Value = "  +2.33%  "
$ws.Cells.Item(49,5).Value = "  +3.64%  "
$ws.Cells.Item(50,5).Value = "  +0.65%  "
$ws.Cells.Item(51,5).Value = "  -0.28%  "

$excel.CutCopyMode = $false

